$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (some use "." as thousands separator, e.g. "28.469.31"),
# so force Text format before assigning to avoid Excel auto-converting to a number,
# then restore the default (Normal) style so formatting matches the source file.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.469.31'
$ws.Range("E2").Value = '  +2.29%  '
$ws.Range("D3").Value = '1.827.60'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '315.62'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.5053'
$ws.Range("E7").Value = '  -5.29%  '
$ws.Range("D8").Value = '0.3916'
$ws.Range("E8").Value = '  +1.58%  '
$ws.Range("D9").Value = '0.07684'
$ws.Range("E9").Value = '  +3.31%  '
$ws.Range("D10").Value = '41.95'
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("E11").Value = '  +2.84%  '
$ws.Range("D12").Value = '21.10'
$ws.Range("E12").Value = '  +4.05%  '
$ws.Range("D13").Value = '6.280'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '7.568'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = '1.824.82'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '93.41'
$ws.Range("E17").Value = '  +6.02%  '
$ws.Range("D18").Value = '0.00001088'
$ws.Range("E18").Value = '  +2.83%  '
$ws.Range("D19").Value = '0.06657'
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  +3.21%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '6.172'
$ws.Range("E22").Value = '  +3.55%  '
$ws.Range("D23").Value = '28.503.38'
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = '2.261'
$ws.Range("E25").Value = '  +8.12%  '
$ws.Range("D26").Value = '157.14'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").Value = '20.64'
$ws.Range("E27").Value = '  +2.68%  '
$ws.Range("D28").Value = '2.037.63'
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("D29").Value = '2.406'
$ws.Range("E29").Value = '  +4.89%  '
$ws.Range("D30").Value = '125.26'
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").Value = '1.132'
$ws.Range("E31").Value = '  +3.35%  '
$ws.Range("D32").Value = '0.1088'
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").Value = '5.679'
$ws.Range("E33").Value = '  +3.47%  '
$ws.Range("D34").Value = '3.665'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '0.07108'
$ws.Range("E35").Value = '  +2.65%  '
$ws.Range("D36").Value = '0.2226'
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '8.918'
$ws.Range("E37").Value = '  +6.47%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02329'
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").Value = '0.6269'
$ws.Range("D41").Value = '11.22'
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("D42").Value = '1.187'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -1.16%  '
$ws.Range("D45").Value = '13.40'
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").Value = '0.5918'
$ws.Range("E46").Value = '  +3.96%  '
$ws.Range("D47").Value = '3.717'
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").Value = '124.63'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = '1.983'
$ws.Range("E49").Value = '  +3.87%  '
$ws.Range("D50").Value = '1.192'
$ws.Range("E50").Value = '  +1.85%  '
$ws.Range("E51").Value = '  +1.84%  '

$ws.Range("D2:D51").Style = "Normal"
